$d = $word.ActiveDocument

# Locate the two paragraphs this edit touches by their content instead of a
# hard-coded index, so the script keeps working even if paragraphs shift.
$addParaIndex = $null
$normalParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($addParaIndex -eq $null -and $txt -like "*Multi threading*") {
        $addParaIndex = $i
    }
    if ($normalParaIndex -eq $null -and $txt -like "*vec3 normal*") {
        $normalParaIndex = $i
    }
}

# --- Change 1 -------------------------------------------------------------
# "Add Multi threading to the current ray tracer I'm building or CUDA it."
# used to be split across three runs (with spell-check proofErr markers
# bracketing "Multi threading"). Collapse it to a single plain run by
# running Find/Replace over the exact same text - Word's replace engine
# writes the replacement back out as one run and drops the proofErr marks.
$addParagraph = $d.Paragraphs.Item($addParaIndex).Range
$addParagraph.Find.Execute(
    "Add Multi threading to the current ray tracer I’m building or CUDA it.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Add Multi threading to the current ray tracer I’m building or CUDA it.",
    2
)

# --- Change 2 -------------------------------------------------------------
# The "vec3 normal = unit_vector(hitPoint - center);" paragraph used to
# carry the (hidden) _GoBack bookmark right before its first run. Remove
# the bookmark from there...
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ...then, right after that paragraph, insert two blank paragraphs followed
# by a paragraph that only contains the _GoBack bookmark (replacing the
# single trailing blank paragraph that used to close out the document).
$normalParagraph = $d.Paragraphs.Item($normalParaIndex)
$insertionPoint = $d.Range($normalParagraph.Range.End, $normalParagraph.Range.End)
$newParasXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' +
               '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' +
               '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$insertionPoint.InsertXML($newParasXml)
